# Applies the "Added button stub and updated tasks" edit to Lab 1 - USART.docx
#
#  - reword the "Hello World string / button" bullet and swap its closing
#    sentence for a note about the empty button-detection stub
#  - tighten the echo-bullet's wording and add a tip about disabling the
#    tx interrupt so the echo is visible
#  - rewrite the command-processor bullet to describe a character-driven
#    command instead of the old digit 0-7 / LED scheme
#  - delete the "For those who like a challenge" (interrupt-driven receive)
#    bullet entirely
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("WARNING: find/replace did not match: " + $find)
    }
}

# ---------------------------------------------------------------------------
# 1) "Hello World string ... button" bullet: reword + replace trailing sentence
# ---------------------------------------------------------------------------
Replace-Text `
    "Modify the firmware so the Hello World string is only transmitted when a button is pressed.  Don’t transmit it continuously.  You may want to look at the Lab0 code to get hints on the button code." `
    "Modify the firmware so the Hello World string is only transmitted when the button is pressed.  Don’t transmit it continuously.  There is an empty section of code for detecting the button in the main loop." | Out-Null

# ---------------------------------------------------------------------------
# 2) "echoed on the terminal" bullet: drop trailing space, add tx-interrupt hint
# ---------------------------------------------------------------------------
Replace-Text `
    "be transmitted back to the PC. " `
    "be transmitted back to the PC.  You will want to disable the tx interrupt first so you can see the echoed text. " | Out-Null

# ---------------------------------------------------------------------------
# 3) "Create a simple command processor" bullet: new command-processor wording
# ---------------------------------------------------------------------------
Replace-Text `
    "Create a simple command processor.  When the PIC receives a digit 0-7 from the PC it should toggle the corresponding LED.  Any other character can be ignored (but should still be echoed.)  You will want to remove the existing lines that toggle LED 0 and 1 on receive." `
    "Create a simple command processor.  When the PIC receives a character of your choice is should execute a simple command.  For example you could set it so an ‘m’ displays a message on the display and an ‘L’ toggles one of the LEDs.  You should have at least two or three commands. You should delete the code that displays the text on the LCD and replace it with your command code." | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the whole "For those who like a challenge..." bullet paragraph
# ---------------------------------------------------------------------------
$deleted = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "For those who like a challenge*") {
        $para.Range.Delete()
        $deleted = $true
        break
    }
}
if (-not $deleted) {
    Write-Output "WARNING: challenge paragraph not found for deletion"
}
